$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Second-count ("Count_Number":"2") transfer-votes column, mirroring the
# existing first-count formulas in column D but reading from column E
# (the per-candidate transfer figures already on the sheet) and bumping
# the running "id" counter by 6 (one extra row) relative to column D.
for ($row = 9; $row -le 14; $row++) {
    $src = $row - 8
    $formula = '="{""Candidate_First_Pref_Votes"":"""&$D' + $src + '&""",' + `
        '""Status"":"""",' + `
        '""Occurred_On_Count"":"""",' + `
        '""Surname"":"""&$B' + $src + '&""",' + `
        '""Firstname"":"""&$A' + $src + '&""",' + `
        '""Constituency_Number"":""2"",' + `
        '""Party_Name"":"""&$C' + $src + '&""",' + `
        '""Candidate_Id"":"""&ROW()&""",' + `
        '""Count_Number"":""2"",' + `
        '""Transfers"":"""&E' + $src + '&""",' + `
        '""id"":"&ROW()-3&",' + `
        '""Total_Votes"":"""&$D' + $src + '+E' + $src + '&"""},"'
    $ws.Range("E$row").Formula = $formula
}

# Match the author's final selection: the new second-count column.
$ws.Range("E9:E14").Select()
